$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (last-changed) date column for rows 2-15
# from 45171 (2023-09-02) to 45172 (2023-09-03), keeping existing formatting.
$ws.Range("C2:C15").Value = 45172
